# Update metadata listing: new file names, new dates, new counts, and a new
# 14th row that did not exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = 43999.37777777778
$ws.Cells.Item(1, 2).Value = 115
$ws.Cells.Item(1, 3).Value = "01._CARATULA.pdf"

# Row 2
$ws.Cells.Item(2, 1).Value = 44064.87149305556
$ws.Cells.Item(2, 2).Value = 4
$ws.Cells.Item(2, 3).Value = "03._aUTO_ADMITE_demaNDA.pdf"

# Row 3
$ws.Cells.Item(3, 1).Value = 44064.87149305556
$ws.Cells.Item(3, 2).Value = 4
$ws.Cells.Item(3, 3).Value = "02._acta_de_reparto.pdf"

# Row 4
$ws.Cells.Item(4, 1).Value = 44123.48692129629
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = "04._MEMORIAL.pdf"

# Row 5
$ws.Cells.Item(5, 1).Value = 44138.73788194444
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = "05._CONSTANCIA_21-03-2021.pdf"

# Row 6
$ws.Cells.Item(6, 1).Value = 44146.6043287037
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = "06._notificación_19.04.2021_DEMANDADO.pdf"

# Row 7
$ws.Cells.Item(7, 1).Value = 44148.70625
$ws.Cells.Item(7, 2).Value = 4
$ws.Cells.Item(7, 3).Value = "07._MEMORIAL,_NO_ACEPTA_DESIGNACION.pdf"

# Row 8
$ws.Cells.Item(8, 1).Value = 44169.64605324074
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = "08._AUTO_NOMBRA_CURADOR.pdf"

# Row 9
$ws.Cells.Item(9, 1).Value = 44181.51538194445
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = "09._acuse_recibido.pdf"

# Row 10
$ws.Cells.Item(10, 1).Value = 44184
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = "10._AcEpTa_dEsIgNaCióN.pdf"

# Row 11
$ws.Cells.Item(11, 1).Value = 44214
$ws.Cells.Item(11, 2).Value = 2
$ws.Cells.Item(11, 3).Value = "11._NOTIFICACION._CURADOR..pdf"

# Row 12
$ws.Cells.Item(12, 1).Value = 44228
$ws.Cells.Item(12, 2).Value = 4
$ws.Cells.Item(12, 3).Value = "12._CONSTESTACION_CURADOR_AD_LITEM.pdf"

# Row 13
$ws.Cells.Item(13, 1).Value = 44231
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = "13._MEMORIAL-SOLICITUD-CELERIDAD.pdf"

# Row 14 (new row)
$ws.Cells.Item(14, 1).Value = 44246.64069444445
$ws.Cells.Item(14, 2).Value = 4
$ws.Cells.Item(14, 3).Value = "14._AUTO-ORDENA--SEGUIR...ADELANTE-EJECUCION.pdf"

# Make sure the new row picks up the same date-time number format as the
# rest of column A.
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat
